$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N10").Value = 177670.86
$ws.Range("O10").Value = 177555.36
$ws.Range("O15").Value = 1684.05
$ws.Range("N21").Value = 720590.33
$ws.Range("N29").Value = 202098
$ws.Range("N30").Value = 18940.7
$ws.Range("N34").Value = 31962.1
